# Corrected salary tax on examples and rerun them
#
# Applies the tax-correction edit to three sheets: "private", "Income" and
# "Folketrygden". For each sheet:
#   - column G ("Utgift") and column Y ("Skattbar") values are corrected
#     for a block of rows (the re-run salary-tax numbers)
#   - column H ("% Endr") gets a 0.5 (50%) value added where it was blank
#   - column G's width grows slightly (it now shows bigger re-computed
#     numbers) - re-applied via AutoFit/ColumnWidth as close as the host
#     allows, since content width changed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "private"
# Rows 38-72: column G and column Y corrected. No H column changes here.
# ---------------------------------------------------------------------
$wsPrivate = $wb.Worksheets.Item("private")

$privateG = @{}
$privateY = @{}
for ($r = 38; $r -le 52; $r++) { $privateG[$r] = 240000.0; $privateY[$r] = -146723.0 }
for ($r = 53; $r -le 56; $r++) { $privateG[$r] = 14900.0;  $privateY[$r] = -371823.0 }
$privateG[57] = 104900.0; $privateY[57] = -281823.0
for ($r = 58; $r -le 67; $r++) { $privateG[$r] = 104900.0; $privateY[$r] = -198700.0 }
for ($r = 68; $r -le 72; $r++) { $privateG[$r] = 92900.0;  $privateY[$r] = -210700.0 }

for ($r = 38; $r -le 72; $r++) {
    $wsPrivate.Cells.Item($r, 7).Value = $privateG[$r]
    $wsPrivate.Cells.Item($r, 25).Value = $privateY[$r]
}

$wsPrivate.Columns.Item(7).ColumnWidth = 8.45

# ---------------------------------------------------------------------
# Sheet 2: "Income"
# Rows 6-37: column H gains a 0.5 value (was blank).
# Rows 38-72: column G, H and Y corrected (H also goes blank -> 0.5).
# ---------------------------------------------------------------------
$wsIncome = $wb.Worksheets.Item("Income")

for ($r = 6; $r -le 37; $r++) {
    $wsIncome.Cells.Item($r, 8).Value = 0.5
}

$incomeG = @{}
$incomeY = @{}
for ($r = 38; $r -le 52; $r++) { $incomeG[$r] = 240000.0; $incomeY[$r] = 60000.0 }
for ($r = 53; $r -le 67; $r++) { $incomeG[$r] = 14900.0;  $incomeY[$r] = -165100.0 }
for ($r = 68; $r -le 72; $r++) { $incomeG[$r] = 2900.0;   $incomeY[$r] = -177100.0 }

for ($r = 38; $r -le 72; $r++) {
    $wsIncome.Cells.Item($r, 7).Value = $incomeG[$r]
    $wsIncome.Cells.Item($r, 8).Value = 0.5
    $wsIncome.Cells.Item($r, 25).Value = $incomeY[$r]
}

$wsIncome.Columns.Item(7).ColumnWidth = 8.45

# ---------------------------------------------------------------------
# Sheet 3: "Folketrygden"
# Rows 6-72: column H gains a 0.5 value (was blank) everywhere.
# Rows 57-72 additionally: column G and Y corrected (0 -> 90000 / 180000 -> 90000).
# ---------------------------------------------------------------------
$wsFolke = $wb.Worksheets.Item("Folketrygden")

for ($r = 6; $r -le 72; $r++) {
    $wsFolke.Cells.Item($r, 8).Value = 0.5
}

for ($r = 57; $r -le 72; $r++) {
    $wsFolke.Cells.Item($r, 7).Value = 90000.0
    $wsFolke.Cells.Item($r, 25).Value = 90000.0
}

$wsFolke.Columns.Item(7).ColumnWidth = 7.25
